$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E hold numeric-looking text (prices with "." thousands
# separators, trailing/leading zeros, percentages with padding spaces)
# that must stay literal strings exactly as scraped. Force Text format
# before writing so Excel does not reinterpret/round the string as a
# number, then restore the "Normal" cell style so no stray number format
# is left behind on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '28.401.36'
Set-TextValue $ws.Range("E2") '  +3.78%  '

Set-TextValue $ws.Range("D3") '1.808.14'
Set-TextValue $ws.Range("E3") '  +1.48%  '

Set-TextValue $ws.Range("D4") '0.9988'
Set-TextValue $ws.Range("E4") '  +0.20%  '

Set-TextValue $ws.Range("D5") '316.09'
Set-TextValue $ws.Range("E5") '  +0.64%  '

Set-TextValue $ws.Range("D6") '0.9989'
Set-TextValue $ws.Range("E6") '  +0.16%  '

Set-TextValue $ws.Range("D7") '0.5512'
Set-TextValue $ws.Range("E7") '  +5.91%  '

Set-TextValue $ws.Range("D8") '0.3861'
Set-TextValue $ws.Range("E8") '  +6.77%  '

Set-TextValue $ws.Range("D9") '0.07596'
Set-TextValue $ws.Range("E9") '  +3.24%  '

Set-TextValue $ws.Range("D10") '42.66'
Set-TextValue $ws.Range("E10") '  +0.19%  '

Set-TextValue $ws.Range("D11") '1.125'
Set-TextValue $ws.Range("E11") '  +2.42%  '

Set-TextValue $ws.Range("D12") '0.9989'
Set-TextValue $ws.Range("E12") '  +0.26%  '

Set-TextValue $ws.Range("D13") '21.16'
Set-TextValue $ws.Range("E13") '  +2.30%  '

Set-TextValue $ws.Range("D14") '6.201'
Set-TextValue $ws.Range("E14") '  +2.01%  '

Set-TextValue $ws.Range("D15") '7.355'
Set-TextValue $ws.Range("E15") '  +5.05%  '

Set-TextValue $ws.Range("D16") '1.804.10'
Set-TextValue $ws.Range("E16") '  +1.87%  '

Set-TextValue $ws.Range("D17") '92.33'
Set-TextValue $ws.Range("E17") '  +4.21%  '

Set-TextValue $ws.Range("D18") '0.00001070'
Set-TextValue $ws.Range("E18") '  +2.18%  '

Set-TextValue $ws.Range("D19") '0.06448'
Set-TextValue $ws.Range("E19") '  +0.25%  '

Set-TextValue $ws.Range("D20") '0.9990'
Set-TextValue $ws.Range("E20") '  +0.16%  '

Set-TextValue $ws.Range("D21") '17.37'
Set-TextValue $ws.Range("E21") '  +3.67%  '

Set-TextValue $ws.Range("E22") '  +2.61%  '

Set-TextValue $ws.Range("D23") '28.410.32'
Set-TextValue $ws.Range("E23") '  +3.61%  '

Set-TextValue $ws.Range("E24") '  +0.96%  '

Set-TextValue $ws.Range("D25") '2.130'
Set-TextValue $ws.Range("E25") '  +2.91%  '

Set-TextValue $ws.Range("D26") '158.44'
Set-TextValue $ws.Range("E26") '  +2.28%  '

Set-TextValue $ws.Range("D27") '20.74'
Set-TextValue $ws.Range("E27") '  +2.57%  '

Set-TextValue $ws.Range("D28") '2.401'
Set-TextValue $ws.Range("E28") '  +1.91%  '

Set-TextValue $ws.Range("D29") '2.012.81'
Set-TextValue $ws.Range("E29") '  +2.05%  '

Set-TextValue $ws.Range("D30") '123.88'
Set-TextValue $ws.Range("E30") '  +1.78%  '

Set-TextValue $ws.Range("D31") '1.128'
Set-TextValue $ws.Range("E31") '  +5.39%  '

Set-TextValue $ws.Range("D32") '0.1020'
Set-TextValue $ws.Range("E32") '  +5.15%  '

Set-TextValue $ws.Range("D33") '5.759'
Set-TextValue $ws.Range("E33") '  +2.92%  '

Set-TextValue $ws.Range("D34") '3.674'
Set-TextValue $ws.Range("E34") '  +2.09%  '

Set-TextValue $ws.Range("D35") '0.2323'
Set-TextValue $ws.Range("E35") '  +14.13%  '

Set-TextValue $ws.Range("D36") '0.06440'
Set-TextValue $ws.Range("E36") '  +7.24%  '

Set-TextValue $ws.Range("D37") '0.02322'
Set-TextValue $ws.Range("E37") '  +3.89%  '

Set-TextValue $ws.Range("D38") '8.838'
Set-TextValue $ws.Range("E38") '  +10.36%  '

Set-TextValue $ws.Range("D39") '11.67'
Set-TextValue $ws.Range("E39") '  +3.79%  '

Set-TextValue $ws.Range("D40") '5.077'
Set-TextValue $ws.Range("E40") '  +4.75%  '

Set-TextValue $ws.Range("D41") '0.6424'
Set-TextValue $ws.Range("E41") '  +4.57%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D42") '1.161'
Set-TextValue $ws.Range("E42") '  +1.54%  '

$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range("D43") '0.9988'
Set-TextValue $ws.Range("E43") '  +0.31%  '

Set-TextValue $ws.Range("D44") '1.380'
Set-TextValue $ws.Range("E44") '  -3.56%  '

Set-TextValue $ws.Range("D45") '13.53'
Set-TextValue $ws.Range("E45") '  +2.06%  '

Set-TextValue $ws.Range("D46") '0.5987'
Set-TextValue $ws.Range("E46") '  +3.64%  '

Set-TextValue $ws.Range("D47") '3.686'
Set-TextValue $ws.Range("E47") '  +1.53%  '

Set-TextValue $ws.Range("D48") '125.23'
Set-TextValue $ws.Range("E48") '  +3.07%  '

Set-TextValue $ws.Range("D49") '1.990'
Set-TextValue $ws.Range("E49") '  +5.11%  '

Set-TextValue $ws.Range("D50") '1.148'
Set-TextValue $ws.Range("E50") '  +3.36%  '

Set-TextValue $ws.Range("D51") '0.06905'
Set-TextValue $ws.Range("E51") '  +2.80%  '
